# Commit: "nueva implementacion de scripts."
# Duplicate the "Login" sheet into a new "LoginDavid" sheet (placed right
# after "Login"), point it at a new david.pena test account, clear out the
# old "new project" verification step (row 8) on both Login sheets, and
# move the active/selected tab from "Objects" to "Login".

$wb = $excel.ActiveWorkbook

$login = $wb.Worksheets.Item("Login")

# --- Create the "LoginDavid" sheet as a copy of "Login", right after it ---
$login.Copy($null, $login)
$loginDavid = $wb.Worksheets.Item(2)
$loginDavid.Name = "LoginDavid"

# --- New credentials for the LoginDavid sheet ---
$loginDavid.Range("D5").Value = "davidpena"
$loginDavid.Range("D6").Value = "Hexaware123"

# --- Clear the old "verify new project" step (row 8) on both sheets ---
foreach ($ws in @($login, $loginDavid)) {
    $ws.Range("A8").ClearContents()
    $ws.Range("B8").ClearContents()
    $ws.Range("C8").ClearContents()
    $ws.Range("D8").ClearContents()
    $ws.Range("E8").ClearContents()
    $ws.Range("F8").ClearContents()
    $ws.Range("I8").ClearContents()
}

# --- Sheet view / selection tweaks ---
$login.Range("C12").Select()
$loginDavid.Range("F11").Select()

$objects = $wb.Worksheets.Item("Objects")
$parameters = $wb.Worksheets.Item("Parameters")

# Move the active tab to "Login" (matches the new tabSelected="1" there)
$login.Activate()
